{"js": "async (context) => {\n  // The template had an accidentally duplicated \"ng\u00e2n h\u00e0ng \" right before\n  // the bank name \"Shinhan\" (i.e. \"...T\u1ea1i Ng\u00e2n h\u00e0ng ng\u00e2n h\u00e0ng Shinhan\n  // Vi\u1ec7t Nam...\"). Remove the stray, lower-case duplicate \"ng\u00e2n h\u00e0ng \"\n  // so the sentence reads \"...T\u1ea1i Ng\u00e2n h\u00e0ng Shinhan Vi\u1ec7t Nam...\".\n  const body = context.document.body;\n\n  const results = body.search(\"ng\u00e2n h\u00e0ng \", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Delete the duplicated phrase in place; the surrounding runs\n    // (\"h\u00e0ng \" before it and the bold \"Shinhan\" after it) are left intact.\n    results.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The template had an accidentally duplicated \"ng\u00e2n h\u00e0ng \" right before\n# the bank name \"Shinhan\" (i.e. \"...T\u1ea1i Ng\u00e2n h\u00e0ng ng\u00e2n h\u00e0ng Shinhan\n# Vi\u1ec7t Nam...\"). Remove the stray, lower-case duplicate \"ng\u00e2n h\u00e0ng \"\n# so the sentence reads \"...T\u1ea1i Ng\u00e2n h\u00e0ng Shinhan Vi\u1ec7t Nam...\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.MatchCase = $true\n$find.Text = \"ng\u00e2n h\u00e0ng \"\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
